# Generate Report for Handback
# Update the timestamp strings that record when the handback report was generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 1be6c6d0... row
$wsOverview.Range("G3").Value = "2016-10-19 10:57:24"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the 1be6c6d0... row
$wsZhCn.Range("H3").Value = "2016-10-19 10:57:13"
$wsZhCn.Range("K3").Value = "2016-10-19 10:57:57"

# de-de sheet: Correspond Handback DateTime for the 1be6c6d0... row
# (Correspond Handoff Datetime on de-de H3 shares the same string as Overview!G3)
$wsDeDe.Range("H3").Value = "2016-10-19 10:57:24"
$wsDeDe.Range("K3").Value = "2016-10-19 10:58:15"
